# Test data added for Italy
#
# 1. Update the selection on the "Slovakia" sheet (it is no longer the
#    active/selected tab once the new sheet is added).
# 2. Add a new worksheet "Itlay" after "Slovakia" by copying the
#    "Germany" sheet (same template layout/styling as the other country
#    sheets) and set its Jira reference cell + selection.

$wb = $excel.ActiveWorkbook

# --- Slovakia: change the remembered selection, it will stop being the
#     active tab once the new sheet is inserted/activated below.
$slovakia = $wb.Worksheets.Item("Slovakia")
$slovakia.Range("C15").Select()

# --- Add the new "Itlay" sheet after "Slovakia", using "Germany" as the
#     template (identical column widths / styles / merged cells).
$germany = $wb.Worksheets.Item("Germany")
$germany.Copy($null, $slovakia)
$italy = $wb.Worksheets.Item($wb.Worksheets.Count)
$italy.Name = "Itlay"

# --- Fill in the Jira/user-story reference for the new market and
#     leave the selection where the author left it.
$italy.Range("B4").Value = "NGC-3145/T2219"
$italy.Range("B4").Select()
